# ZooBC - Feature List.xlsx
# Swap the "% Complete" values for the two "Product Reviews" features:
#   - "Display product review on product detail"            (row 30): 0 -> 0.5
#   - "Restrict to one review per user on item purchased"    (row 31): 0.5 -> 0
# The pivot table ("Pivot Table Report") reflects these same two features
# (rows 8-9 on the pivot side) and recalculates automatically.
# Finally, move the active selection to K19 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D30").Value = 0.5
$ws.Range("D31").Value = 0

$ws.Range("K19").Select()
